$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the sheet's original (visible) gridlines setting intact while we work.
$excel.ActiveWindow.DisplayGridlines = $true

# Finished going through Ginger's comments: fill in the remaining replies
# in column L (mapper response) for the last few rows of the table.
$ws.Range("L18").Value = "Done."
$ws.Range("L19").Value = "Done."
$ws.Range("L20").Value = "Done. I'm curious, how did you catch this one? It's quite close to the intersection of road and stream, but not exactly."
$ws.Range("L21").Value = "OK."

# Scroll so row 17 is at the top and leave the last edited cell selected,
# matching where the editor ended up after finishing these replies.
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L21").Select() | Out-Null
